# Atalanta 2020 matches - fill in results for 5 more fixtures
# (Atalanta-Verona, Atalanta-Fiorentina, Juventus-Atalanta, Atalanta-Roma,
#  Bologna-Atalanta) which previously only had the fixture (home/away teams)
# without xG / goals figures, and shift the remaining still-unplayed
# fixtures list down so it keeps starting right after the last completed
# match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the five newly completed fixtures (row 10 already had the
#     correct home/away teams, rows 11-14 need their home/away swapped to
#     the next fixtures on the calendar) ---

# Row 10: Atalanta 0-2 Verona (xG 1.20139 - 1.78911)
$ws.Range("D10:G10").NumberFormat = "@"
$ws.Range("D10").Value = "1.20139"
$ws.Range("E10").Value = "1.78911"
$ws.Range("F10").Value = "0"
$ws.Range("G10").Value = "2"

# Row 11: Atalanta 3-0 Fiorentina (xG 2.8038 - 0.461794)
$ws.Range("B11").Value = "Atalanta"
$ws.Range("C11").Value = "Fiorentina"
$ws.Range("D11:G11").NumberFormat = "@"
$ws.Range("D11").Value = "2.8038"
$ws.Range("E11").Value = "0.461794"
$ws.Range("F11").Value = "3"
$ws.Range("G11").Value = "0"

# Row 12: Juventus 1-1 Atalanta (xG 2.71039 - 0.902039)
$ws.Range("B12").Value = "Juventus"
$ws.Range("C12").Value = "Atalanta"
$ws.Range("D12:G12").NumberFormat = "@"
$ws.Range("D12").Value = "2.71039"
$ws.Range("E12").Value = "0.902039"
$ws.Range("F12").Value = "1"
$ws.Range("G12").Value = "1"

# Row 13: Atalanta 4-1 Roma (xG 1.91426 - 0.712758)
$ws.Range("B13").Value = "Atalanta"
$ws.Range("C13").Value = "Roma"
$ws.Range("D13:G13").NumberFormat = "@"
$ws.Range("D13").Value = "1.91426"
$ws.Range("E13").Value = "0.712758"
$ws.Range("F13").Value = "4"
$ws.Range("G13").Value = "1"

# Row 14: Bologna 2-2 Atalanta (xG 1.08799 - 1.89381)
$ws.Range("B14").Value = "Bologna"
$ws.Range("C14").Value = "Atalanta"
$ws.Range("D14:G14").NumberFormat = "@"
$ws.Range("D14").Value = "1.08799"
$ws.Range("E14").Value = "1.89381"
$ws.Range("F14").Value = "2"
$ws.Range("G14").Value = "2"

# --- Shift the remaining still-unplayed fixtures down by one slot so the
#     schedule continues right after the match played in row 14 ---

$ws.Range("B15").Value = "Atalanta"
$ws.Range("C15").Value = "Sassuolo"

$ws.Range("C16").Value = "Parma Calcio 1913"

$ws.Range("B17").Value = "Benevento"
$ws.Range("C17").Value = "Atalanta"

$ws.Range("B18").Value = "Atalanta"
$ws.Range("C18").Value = "Genoa"

$ws.Range("B19").Value = "Udinese"
$ws.Range("C19").Value = "Atalanta"
